$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New flashcard row (row 38): ID / Question / Answer.
$ws.Range("A38").Value = "2023-10-17 13:42:35 6 question_6_2"
$ws.Range("B38").Value = "6 question"

# The answer "6" looks numeric, but the source data stores it as text
# (same as every other cell in this sheet). Force text entry by
# flipping the cell to a text format before assigning the value, then
# restore the original "General" format so the cell keeps the sheet's
# normal style.
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "6"
$ws.Range("C38").NumberFormat = "General"
